# Fix typo in 3rd part of max-cut/min-flow proof:
# "Proof by contradiction, using weak duality property"
#   -> "Proof by contradiction"
# (slide "C -> A: if f is max-flow, then no augmenting path")

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetParaIndex = -1

foreach ($s in $p.Slides) {
    foreach ($shape in $s.Shapes) {
        if (-not $shape.HasTextFrame) { continue }
        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($i = 1; $i -le $paraCount; $i++) {
            $para = $tr.Paragraphs($i, 1)
            $paraText = $para.Text.TrimEnd("`r")
            if ($paraText -eq "Proof by contradiction, using weak duality property") {
                $targetSlide = $s
                $targetShape = $shape
                $targetParaIndex = $i
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the 'Proof by contradiction, using weak duality property' paragraph."
}

$tr = $targetShape.TextFrame.TextRange
$para = $tr.Paragraphs($targetParaIndex, 1)
$paraStart = $para.Start

# Replace the paragraph text, dropping ", using weak duality property".
$para.Text = "Proof by contradiction"

# Re-apply the "by contradiction" portion as its own run so it is stored
# as a distinct <a:r> from the "Proof " run (matches the authored split).
$secondRun = $tr.Characters($paraStart + 6, 16)
$secondRun.Text = "by contradiction"
